$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete early-year rows (2000年, 2005年, 2008年, 2009年).
# Deleting rows 2:5 shifts the 2010年..2020年 rows (previously rows 6..16) up to rows 2..12.
$ws.Rows("2:5").Delete()

# Append a new row for 2021年 at the bottom (row 13), copying the formatting
# (bold/border/centered style) from the row above (row 12, 2020年).
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "2021年"
$ws.Range("U13").Value = 7093
